$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column stores values as plain text (e.g. "43.22", "68.224.66").
# Force affected cells to Text format first so Excel does not silently
# reinterpret/round numeric-looking strings as floating point numbers.
foreach ($r in @(2,3,5,6,7,11,12,13,14,15,16,18,20,21,22,23,24,25,26,27,28,30,31,32,33,34,35,36,38,39,40,41,42,43,44,45,46,48,51)) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
}

$ws.Cells.Item(2, 4).Value = '68.224.66'
$ws.Cells.Item(2, 5).Value = '  +1.01%  '

$ws.Cells.Item(3, 4).Value = '3.906.85'
$ws.Cells.Item(3, 5).Value = '  -0.81%  '

$ws.Cells.Item(4, 5).Value = '  +0.17%  '

$ws.Cells.Item(5, 4).Value = '486.16'
$ws.Cells.Item(5, 5).Value = '  +2.62%  '

$ws.Cells.Item(6, 4).Value = '146.16'
$ws.Cells.Item(6, 5).Value = '  -2.39%  '

$ws.Cells.Item(7, 4).Value = '0.622'
$ws.Cells.Item(7, 5).Value = '  -0.68%  '

$ws.Cells.Item(8, 5).Value = '  -0.10%  '

$ws.Cells.Item(10, 5).Value = '  -0.48%  '

$ws.Cells.Item(11, 4).Value = '0.0000341'
$ws.Cells.Item(11, 5).Value = '  -1.66%  '

$ws.Cells.Item(12, 4).Value = '43.22'
$ws.Cells.Item(12, 5).Value = '  -0.52%  '

$ws.Cells.Item(13, 4).Value = '10.85'
$ws.Cells.Item(13, 5).Value = '  +4.11%  '

$ws.Cells.Item(14, 4).Value = '4.534.57'
$ws.Cells.Item(14, 5).Value = '  -0.81%  '

$ws.Cells.Item(15, 4).Value = '3.900.89'
$ws.Cells.Item(15, 5).Value = '  -1.33%  '

$ws.Cells.Item(16, 4).Value = '14.19'
$ws.Cells.Item(16, 5).Value = '  -5.13%  '

$ws.Cells.Item(18, 4).Value = '19.97'
$ws.Cells.Item(18, 5).Value = '  +0.09%  '

$ws.Cells.Item(19, 5).Value = '  -1.42%  '

$ws.Cells.Item(20, 4).Value = '68.331.83'
$ws.Cells.Item(20, 5).Value = '  +0.87%  '

$ws.Cells.Item(21, 4).Value = '432.74'
$ws.Cells.Item(21, 5).Value = '  -0.52%  '

$ws.Cells.Item(22, 4).Value = '3.50'
$ws.Cells.Item(22, 5).Value = '  +3.24%  '

$ws.Cells.Item(23, 4).Value = '14.91'
$ws.Cells.Item(23, 5).Value = '  +3.30%  '

$ws.Cells.Item(24, 4).Value = '87.89'
$ws.Cells.Item(24, 5).Value = '  +0.34%  '

$ws.Cells.Item(25, 4).Value = '11.19'
$ws.Cells.Item(25, 5).Value = '  +13.36%  '

$ws.Cells.Item(26, 4).Value = '11.24'
$ws.Cells.Item(26, 5).Value = '  +10.41%  '

$ws.Cells.Item(27, 4).Value = '3.58'
$ws.Cells.Item(27, 5).Value = '  -1.06%  '

$ws.Cells.Item(28, 4).Value = '37.98'
$ws.Cells.Item(28, 5).Value = '  -1.63%  '

$ws.Cells.Item(29, 5).Value = '  +0.86%  '

$ws.Cells.Item(30, 4).Value = '720.83'
$ws.Cells.Item(30, 5).Value = '  +0.21%  '

$ws.Cells.Item(31, 4).Value = '13.69'
$ws.Cells.Item(31, 5).Value = '  +2.13%  '

$ws.Cells.Item(32, 4).Value = '0.129'
$ws.Cells.Item(32, 5).Value = '  -2.27%  '

$ws.Cells.Item(33, 4).Value = '2.92'
$ws.Cells.Item(33, 5).Value = '  +4.15%  '

$ws.Cells.Item(34, 4).Value = '6.18'
$ws.Cells.Item(34, 5).Value = '  +14.62%  '

$ws.Cells.Item(35, 4).Value = '41.38'
$ws.Cells.Item(35, 5).Value = '  -1.90%  '

$ws.Cells.Item(36, 4).Value = '0.0₃0863'
$ws.Cells.Item(36, 5).Value = '  +4.42%  '

$ws.Cells.Item(37, 5).Value = '  +4.10%  '

$ws.Cells.Item(38, 4).Value = '0.147'
$ws.Cells.Item(38, 5).Value = '  -3.06%  '

$ws.Cells.Item(41, 4).Value = '0.0485'
$ws.Cells.Item(41, 5).Value = '  +1.89%  '

$ws.Cells.Item(42, 4).Value = '2.88'
$ws.Cells.Item(42, 5).Value = '  +12.90%  '

$ws.Cells.Item(43, 4).Value = '3.13'
$ws.Cells.Item(43, 5).Value = '  +2.85%  '

$ws.Cells.Item(44, 4).Value = '2.95'
$ws.Cells.Item(44, 5).Value = '  +4.39%  '

$ws.Cells.Item(47, 5).Value = '  +0.05%  '

$ws.Cells.Item(48, 4).Value = '0.0₆0352'
$ws.Cells.Item(48, 5).Value = '  +38.89%  '

$ws.Cells.Item(49, 5).Value = '  -2.62%  '

$ws.Cells.Item(50, 5).Value = '  -3.93%  '

$ws.Cells.Item(51, 4).Value = '145.05'
$ws.Cells.Item(51, 5).Value = '  -1.70%  '

# Rows 39/40 and 45/46: coin entries were reordered (swapped) in the ranking.
$ws.Cells.Item(39, 2).Value = 'Dai'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Cells.Item(39, 4).Value = '0.997'
$ws.Cells.Item(39, 5).Value = '  -0.12%  '

$ws.Cells.Item(40, 2).Value = 'TheGraph'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Cells.Item(40, 4).Value = '0.391'
$ws.Cells.Item(40, 5).Value = '  +16.15%  '

$ws.Cells.Item(45, 2).Value = 'ApeXProtocol'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Cells.Item(45, 4).Value = '3.46'
$ws.Cells.Item(45, 5).Value = '  +6.89%  '

$ws.Cells.Item(46, 2).Value = 'Stellar'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(46, 4).Value = '0.141'
$ws.Cells.Item(46, 5).Value = '  -0.98%  '
